$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33438
$ws.Range("J3").Value = 33438
$ws.Range("L3").Value = 33438
$ws.Range("N3").Value = -33666
$ws.Range("H76").Value = 4508.154
$ws.Range("I76").Value = 3966.6667
$ws.Range("J76").Value = 4670.6
$ws.Range("K76").Value = 3966.6667
$ws.Range("L76").Value = 4670.6
$ws.Range("M76").Value = -3651.6667
$ws.Range("N76").Value = -5300.6
$ws.Range("H79").Value = 4508.154
$ws.Range("I79").Value = 3966.6667
$ws.Range("J79").Value = 4670.6
$ws.Range("K79").Value = 3966.6667
$ws.Range("L79").Value = 4670.6
$ws.Range("M79").Value = -2874.6667
$ws.Range("N79").Value = -6854.6
$ws.Range("H92").Value = 1179.8
$ws.Range("I92").Value = 1371.2858
$ws.Range("J92").Value = 733
$ws.Range("K92").Value = 1371.2858
$ws.Range("L92").Value = 733
$ws.Range("M92").Value = -123.2858000000001
$ws.Range("N92").Value = -3229
$ws.Range("H98").Value = 1869.1
$ws.Range("I98").Value = 1742.7778
$ws.Range("K98").Value = 1742.7778
$ws.Range("M98").Value = -244.7778000000001
$ws.Range("H102").Value = 33438
$ws.Range("J102").Value = 33438
$ws.Range("L102").Value = 33438
$ws.Range("N102").Value = -39928
$ws.Range("H106").Value = 4000
$ws.Range("I106").Value = 4000
$ws.Range("K106").Value = 4000
$ws.Range("M106").Value = -3369
$ws.Range("H111").Value = 10786.066
$ws.Range("J111").Value = 3700.5
$ws.Range("L111").Value = 11101.5
$ws.Range("N111").Value = -17235.5
$ws.Range("H112").Value = 1101.619
$ws.Range("J112").Value = 1131.7
$ws.Range("L112").Value = 3395.1
$ws.Range("N112").Value = -5611.1
$ws.Range("H122").Value = 1869.1
$ws.Range("I122").Value = 1742.7778
$ws.Range("K122").Value = 5228.3334
$ws.Range("M122").Value = -2778.3334
$ws.Range("H137").Value = 1309.9756
$ws.Range("I137").Value = 1101.0605
$ws.Range("J137").Value = 2171.75
$ws.Range("K137").Value = 3303.1815
$ws.Range("L137").Value = 6515.25
$ws.Range("M137").Value = -753.1815000000001
$ws.Range("N137").Value = -11615.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20005.22
$ws.Range("I32").Value = 3224.622
$ws.Range("J32").Value = 172895.11
$ws.Range("K32").Value = 3224.622
$ws.Range("L32").Value = 172895.11
$ws.Range("M32").Value = -2937.622
$ws.Range("N32").Value = -173469.11
$ws.Range("H122").Value = 1916.9
$ws.Range("I122").Value = 2021.5
$ws.Range("J122").Value = 1498.5
$ws.Range("K122").Value = 6064.5
$ws.Range("L122").Value = 4495.5
$ws.Range("M122").Value = -3614.5
$ws.Range("N122").Value = -9395.5
$ws.Range("H132").Value = 2021.093
$ws.Range("I132").Value = 1862.4054
$ws.Range("J132").Value = 2999.6667
$ws.Range("K132").Value = 5587.216200000001
$ws.Range("L132").Value = 8999.000100000001
$ws.Range("M132").Value = -3057.216200000001
$ws.Range("N132").Value = -14059.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 753.8
$ws.Range("I22").Value = 317
$ws.Range("J22").Value = 2501
$ws.Range("K22").Value = 317
$ws.Range("L22").Value = 2501
$ws.Range("M22").Value = 33
$ws.Range("N22").Value = -3201
$ws.Range("H31").Value = 29248.21
$ws.Range("I31").Value = 1063.3529
$ws.Range("K31").Value = 1063.3529
$ws.Range("M31").Value = -768.3529000000001
$ws.Range("H34").Value = 29248.21
$ws.Range("I34").Value = 1063.3529
$ws.Range("K34").Value = 1063.3529
$ws.Range("M34").Value = -861.3529000000001
$ws.Range("H58").Value = 5102.8716
$ws.Range("I58").Value = 796.96
$ws.Range("J58").Value = 12792
$ws.Range("K58").Value = 796.96
$ws.Range("L58").Value = 12792
$ws.Range("M58").Value = -593.96
$ws.Range("N58").Value = -13198
$ws.Range("H99").Value = 10139.934
$ws.Range("I99").Value = 3315.1667
$ws.Range("J99").Value = 14689.777
$ws.Range("K99").Value = 3315.1667
$ws.Range("L99").Value = 14689.777
$ws.Range("M99").Value = -1817.1667
$ws.Range("N99").Value = -17685.777
$ws.Range("H122").Value = 499
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 498
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 1494
$ws.Range("M122").Value = 950
$ws.Range("N122").Value = -6394
$ws.Range("H126").Value = 10139.934
$ws.Range("I126").Value = 3315.1667
$ws.Range("J126").Value = 14689.777
$ws.Range("K126").Value = 9945.500100000001
$ws.Range("L126").Value = 44069.331
$ws.Range("M126").Value = -7475.500100000001
$ws.Range("N126").Value = -49009.331
$ws.Range("H136").Value = 5102.8716
$ws.Range("I136").Value = 796.96
$ws.Range("J136").Value = 12792
$ws.Range("K136").Value = 2390.88
$ws.Range("L136").Value = 38376
$ws.Range("M136").Value = 159.1199999999999
$ws.Range("N136").Value = -43476
$ws.Range("H137").Value = 34656.5
$ws.Range("I137").Value = 20709
$ws.Range("J137").Value = 37446
$ws.Range("K137").Value = 20709
$ws.Range("L137").Value = 37446
$ws.Range("M137").Value = -15609
$ws.Range("N137").Value = -47646

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 599.2857
$ws.Range("J34").Value = 956.25
$ws.Range("L34").Value = 2868.75
$ws.Range("N34").Value = -3036.75
$ws.Range("H131").Value = 759.48
$ws.Range("I131").Value = 360.0625
$ws.Range("J131").Value = 835.5595
$ws.Range("K131").Value = 1080.1875
$ws.Range("L131").Value = 2506.6785
$ws.Range("M131").Value = 3959.8125
$ws.Range("N131").Value = -12586.6785

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 38000
$ws.Range("J101").Value = 38000
$ws.Range("L101").Value = 38000
$ws.Range("N101").Value = -44490
$ws.Range("H105").Value = 45952.5
$ws.Range("J105").Value = 45952.5
$ws.Range("L105").Value = 45952.5
$ws.Range("N105").Value = -52940.5
$ws.Range("H132").Value = 2960.0454
$ws.Range("I132").Value = 2806.7778
$ws.Range("J132").Value = 3649.75
$ws.Range("K132").Value = 8420.3334
$ws.Range("L132").Value = 10949.25
$ws.Range("M132").Value = -5890.3334
$ws.Range("N132").Value = -16009.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 45193
$ws.Range("J118").Value = 45193
$ws.Range("L118").Value = 45193
$ws.Range("N118").Value = -48507

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 41417.5
$ws.Range("J117").Value = 41417.5
$ws.Range("L117").Value = 41417.5
$ws.Range("N117").Value = -50595.5
$ws.Range("H122").Value = 1615
$ws.Range("I122").Value = 1607.2222
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 4821.6666
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -2371.6666
$ws.Range("N122").Value = -9850
$ws.Range("H126").Value = 1728.4286
$ws.Range("I126").Value = 1518
$ws.Range("K126").Value = 4554
$ws.Range("M126").Value = -2084
